$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 19 (2025Q1) metrics per corrected figures
$ws.Range("C19").Value = 298
$ws.Range("D19").Value = 253
$ws.Range("E19").Value = 45
$ws.Range("F19").Value = 79.31034482758621
